# Re-derive the UTM Easting/Northing (columns B/C), forcing UTM zone 1
# (columns D/E), for every data row on Sheet1, using the existing decimal
# latitude/longitude already present in columns I ("Y") and J ("X").
#
# This mirrors the standard WGS84 lat/lon -> UTM forward transform (the
# same formulation used by the widely used python "utm" package), but is
# evaluated through native Excel worksheet formulas (SIN/COS/TAN/SQRT/PI)
# so the arithmetic/rounding matches Excel's own engine bit for bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Scratch cell, far away from the real data, used to evaluate each
# formula via the live Excel calculation engine and pull back a plain
# numeric value.
$scratch = $ws.Range("ZZ1")

for ($row = 2; $row -le $lastRow; $row++) {

    $lat = $ws.Cells.Item($row, 9).Value()   # column I = Y (decimal latitude)
    $lon = $ws.Cells.Item($row, 10).Value()  # column J = X (decimal longitude)

    if ($lat -eq $null -or $lon -eq $null) {
        continue
    }

    $zone = 1
    $centralLonDeg = (($zone - 1) * 6 - 180 + 3)

    # Full (17-significant-digit) literals so nothing is lost when the
    # numbers are spliced into the formula text.
    $latLit = $lat.ToString("G17")
    $lonLit = $lon.ToString("G17")

    # Shared LET() preamble for both Easting and Northing. Variable names
    # are deliberately non-cell-like (no bare "<letter><digit>" tokens
    # such as e2/n1/c1/a1/t1/m1 - those parse as cell references instead
    # of LET names).
    $pre = "latRad,$latLit*(PI()/180)," +
           "lonRad,$lonLit*(PI()/180)," +
           "clRad,$centralLonDeg*(PI()/180)," +
           "eccv,0.00669438," +
           "epv,eccv/(1-eccv)," +
           "nv,6378137/SQRT(1-eccv*SIN(latRad)^2)," +
           "cv,epv*COS(latRad)^2," +
           "av,COS(latRad)*(lonRad-clRad)," +
           "tv,TAN(latRad)," +
           "mv,6378137*((1-eccv/4-3*eccv^2/64-5*eccv^3/256)*latRad-(3*eccv/8+3*eccv^2/32+45*eccv^3/1024)*SIN(2*latRad)+(15*eccv^2/256+45*eccv^3/1024)*SIN(4*latRad)-(35*eccv^3/3072)*SIN(6*latRad)),"

    $eastingFormula = "=LET($pre 0.9996*nv*(av+av^3/6*(1-tv^2+cv)+av^5/120*(5-18*tv^2+tv^4+72*cv-58*epv))+500000)"
    $northingFormula = "=LET($pre 0.9996*(mv+nv*tv*(av^2/2+av^4/24*(5-tv^2+9*cv+4*cv^2)+av^6/720*(61-58*tv^2+tv^4+600*cv-330*epv))))"

    $scratch.Formula = $eastingFormula
    $easting = $scratch.Value()

    $scratch.Formula = $northingFormula
    $northing = $scratch.Value()

    $ws.Cells.Item($row, 2).Value = $easting    # column B = Y_UTM (easting)
    $ws.Cells.Item($row, 3).Value = $northing   # column C = X_UTM (northing)
    $ws.Cells.Item($row, 4).Value = $zone       # column D = ZoneNumber
    $ws.Cells.Item($row, 5).Value = "W"         # column E = ZoneLetter
}

# Clean up the scratch cell / formula so it doesn't leak into the saved
# workbook.
$scratch.Formula = ""
$scratch.Value = ""
